# Generate Report for Handback
# Re-sorts/fills in the localization-status workbook so that all three
# source files (31837146-..., 4e2f794a-..., e6052173-...) show up as
# "Handed back: in sync with en-US" with full handoff/handback detail,
# ordered alphabetically by file name across all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Hyperlinks.Delete()

$ov.Range("A1").Value2 = "File Name"
$ov.Range("B1").Value2 = "zh-cn"
$ov.Range("C1").Value2 = "de-de"

$ov.Range("A2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.md"
$ov.Range("B2").Value2 = "Handed back: in sync with en-US"
$ov.Range("C2").Value2 = "Handed back: in sync with en-US"

$ov.Range("A3").Value2 = "4e2f794a-6238-4484-ae55-20073f44aa94.md"
$ov.Range("B3").Value2 = "Handed back: in sync with en-US"
$ov.Range("C3").Value2 = "Handed back: in sync with en-US"

$ov.Range("A4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.md"
$ov.Range("B4").Value2 = "Handed back: in sync with en-US"
$ov.Range("C4").Value2 = "Handed back: in sync with en-US"

$ov.Range("A5").Value2 = ".localization-config"
$ov.Range("B5").Value2 = "Not to be localized"
$ov.Range("C5").Value2 = "Not to be localized"

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/31837146-42ff-4698-bfdd-003adfeb227b.md", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/4e2f794a-6238-4484-ae55-20073f44aa94.md", "", "", "4e2f794a-6238-4484-ae55-20073f44aa94.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3c25566890736fe06452d0dea93bfc1800087a25/e2e/e6052173-41f7-4b7d-8df8-c8821d2e3359.md", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Hyperlinks.Delete()

$zh.Range("A1").Value2 = "Source File Name"
$zh.Range("B1").Value2 = "Status"
$zh.Range("C1").Value2 = "Latest Handoff File"
$zh.Range("D1").Value2 = "Latest Handoff Datetime"
$zh.Range("E1").Value2 = "Latest Target File"
$zh.Range("F1").Value2 = "Latest Handback File"
$zh.Range("G1").Value2 = "Latest Handback DateTime"
$zh.Range("H1").Value2 = "Handoff Reason"
$zh.Range("I1").Value2 = "Dependency From"

# Row 2: 31837146-...
$zh.Range("A2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.md"
$zh.Range("B2").Value2 = "Handed back: in sync with en-US"
$zh.Range("C2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf"
$zh.Range("D2").Value2 = "2016-02-24 07:13:35"
$zh.Range("E2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.md"
$zh.Range("F2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf"
$zh.Range("G2").Value2 = "2016-02-24 07:14:22"
$zh.Range("H2").Value2 = "Include"

# Row 3: 4e2f794a-... (mirrors the 31837146 handoff/handback info - dependency)
$zh.Range("A3").Value2 = "4e2f794a-6238-4484-ae55-20073f44aa94.md"
$zh.Range("B3").Value2 = "Handed back: in sync with en-US"
$zh.Range("C3").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf"
$zh.Range("D3").Value2 = "2016-02-24 07:13:35"
$zh.Range("E3").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.md"
$zh.Range("F3").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf"
$zh.Range("G3").Value2 = "2016-02-24 07:14:22"
$zh.Range("H3").Value2 = "Include"

# Row 4: e6052173-...
$zh.Range("A4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.md"
$zh.Range("B4").Value2 = "Handed back: in sync with en-US"
$zh.Range("C4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.zh-cn.xlf"
$zh.Range("D4").Value2 = "2016-02-24 07:11:44"
$zh.Range("E4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.md"
$zh.Range("F4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.zh-cn.xlf"
$zh.Range("G4").Value2 = "2016-02-24 07:12:30"
$zh.Range("H4").Value2 = "Include"

# Row 5: .localization-config
$zh.Range("A5").Value2 = ".localization-config"
$zh.Range("B5").Value2 = "Not to be localized"
$zh.Range("D5").Value2 = "0001-01-01 00:00:00"
$zh.Range("G5").Value2 = "0001-01-01 00:00:00"
$zh.Range("H5").Value2 = "Ignored"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/31837146-42ff-4698-bfdd-003adfeb227b.md", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b05be107cdb740665a450a4fa5948226616442b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/31837146-42ff-4698-bfdd-003adfeb227b.md", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a39754823364991a6935e7d38491de158e9f9465/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/4e2f794a-6238-4484-ae55-20073f44aa94.md", "", "", "4e2f794a-6238-4484-ae55-20073f44aa94.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b05be107cdb740665a450a4fa5948226616442b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/31837146-42ff-4698-bfdd-003adfeb227b.md", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a39754823364991a6935e7d38491de158e9f9465/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3c25566890736fe06452d0dea93bfc1800087a25/e2e/e6052173-41f7-4b7d-8df8-c8821d2e3359.md", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/159d09de854942556170ad791a0eba2902601371/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.zh-cn.xlf", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/598b88f307af80256802d5d76f63d3cb211c603a/e2e/e6052173-41f7-4b7d-8df8-c8821d2e3359.md", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a39754823364991a6935e7d38491de158e9f9465/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.zh-cn.xlf", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Hyperlinks.Delete()

$de.Range("A1").Value2 = "Source File Name"
$de.Range("B1").Value2 = "Status"
$de.Range("C1").Value2 = "Latest Handoff File"
$de.Range("D1").Value2 = "Latest Handoff Datetime"
$de.Range("E1").Value2 = "Latest Target File"
$de.Range("F1").Value2 = "Latest Handback File"
$de.Range("G1").Value2 = "Latest Handback DateTime"
$de.Range("H1").Value2 = "Handoff Reason"
$de.Range("I1").Value2 = "Dependency From"

# Row 2: 31837146-...
$de.Range("A2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.md"
$de.Range("B2").Value2 = "Handed back: in sync with en-US"
$de.Range("C2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf"
$de.Range("D2").Value2 = "2016-02-24 07:13:47"
$de.Range("E2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.md"
$de.Range("F2").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf"
$de.Range("G2").Value2 = "2016-02-24 07:14:43"
$de.Range("H2").Value2 = "Include"

# Row 3: 4e2f794a-... (mirrors the 31837146 handoff/handback info - dependency)
$de.Range("A3").Value2 = "4e2f794a-6238-4484-ae55-20073f44aa94.md"
$de.Range("B3").Value2 = "Handed back: in sync with en-US"
$de.Range("C3").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf"
$de.Range("D3").Value2 = "2016-02-24 07:13:47"
$de.Range("E3").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.md"
$de.Range("F3").Value2 = "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf"
$de.Range("G3").Value2 = "2016-02-24 07:14:43"
$de.Range("H3").Value2 = "Include"

# Row 4: e6052173-...
$de.Range("A4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.md"
$de.Range("B4").Value2 = "Handed back: in sync with en-US"
$de.Range("C4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.de-de.xlf"
$de.Range("D4").Value2 = "2016-02-24 07:11:56"
$de.Range("E4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.md"
$de.Range("F4").Value2 = "e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.de-de.xlf"
$de.Range("G4").Value2 = "2016-02-24 07:12:51"
$de.Range("H4").Value2 = "Include"

# Row 5: .localization-config
$de.Range("A5").Value2 = ".localization-config"
$de.Range("B5").Value2 = "Not to be localized"
$de.Range("D5").Value2 = "0001-01-01 00:00:00"
$de.Range("G5").Value2 = "0001-01-01 00:00:00"
$de.Range("H5").Value2 = "Ignored"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/31837146-42ff-4698-bfdd-003adfeb227b.md", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4adda398144a07be9ae40715f3a22d5b7269d2d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/31837146-42ff-4698-bfdd-003adfeb227b.md", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4b19065c231059748947ca69edefb3dc23967e8e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/4e2f794a-6238-4484-ae55-20073f44aa94.md", "", "", "4e2f794a-6238-4484-ae55-20073f44aa94.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c4adda398144a07be9ae40715f3a22d5b7269d2d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/e2e/31837146-42ff-4698-bfdd-003adfeb227b.md", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4b19065c231059748947ca69edefb3dc23967e8e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf", "", "", "31837146-42ff-4698-bfdd-003adfeb227b.047c5ec904881f5c4d9ade616a56822a922e4c8c.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3c25566890736fe06452d0dea93bfc1800087a25/e2e/e6052173-41f7-4b7d-8df8-c8821d2e3359.md", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ce5a93f9a9f847ed1233217f59e9862ad768698/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.de-de.xlf", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1c3c30ee4503fc0a4f8496aa6bf02ca68eca4f9b/e2e/e6052173-41f7-4b7d-8df8-c8821d2e3359.md", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4b19065c231059748947ca69edefb3dc23967e8e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.de-de.xlf", "", "", "e6052173-41f7-4b7d-8df8-c8821d2e3359.1e4e3286c2b3b6c1c33520f07aa843ea0ab86005.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/186b7e5dc65b887f7e5aa79b88f491b2de99d783/.localization-config", "", "", ".localization-config") | Out-Null
